# Updates the "cryptos" price table on Sheet1 with refreshed values scraped
# on Sun Dec 17 18:34:39 UTC 2023 (see commit message). Each entry below is
# one changed cell: Price (column D), Volume/1h (column E), plus the
# Coin/Link (columns B/C) for the two rows whose ranking swapped places
# (THORChain <-> MultiversX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '42.154.83' },
    @{ Cell = 'E2'; Value = '  -0.82%  ' },
    @{ Cell = 'D3'; Value = '2.224.84' },
    @{ Cell = 'E3'; Value = '  -0.87%  ' },
    @{ Cell = 'E4'; Value = '  -0.23%  ' },
    @{ Cell = 'D5'; Value = '243.13' },
    @{ Cell = 'E5'; Value = '  -0.68%  ' },
    @{ Cell = 'D6'; Value = '0.627' },
    @{ Cell = 'E6'; Value = '  -0.39%  ' },
    @{ Cell = 'D7'; Value = '74.31' },
    @{ Cell = 'E7'; Value = '  -1.73%  ' },
    @{ Cell = 'E8'; Value = '  -0.01%  ' },
    @{ Cell = 'D9'; Value = '0.604' },
    @{ Cell = 'E9'; Value = '  -2.90%  ' },
    @{ Cell = 'D10'; Value = '42.74' },
    @{ Cell = 'E10'; Value = '  -2.47%  ' },
    @{ Cell = 'E11'; Value = '  +1.35%  ' },
    @{ Cell = 'D12'; Value = '6.97' },
    @{ Cell = 'E12'; Value = '  -3.31%  ' },
    @{ Cell = 'E13'; Value = '  +0.25%  ' },
    @{ Cell = 'D14'; Value = '2.557.79' },
    @{ Cell = 'E14'; Value = '  -0.77%  ' },
    @{ Cell = 'D15'; Value = '14.31' },
    @{ Cell = 'E15'; Value = '  -1.82%  ' },
    @{ Cell = 'D16'; Value = '0.837' },
    @{ Cell = 'E16'; Value = '  -2.58%  ' },
    @{ Cell = 'D17'; Value = '2.214.56' },
    @{ Cell = 'E17'; Value = '  -0.68%  ' },
    @{ Cell = 'D18'; Value = '41.956.99' },
    @{ Cell = 'E18'; Value = '  -0.78%  ' },
    @{ Cell = 'E19'; Value = '  +4.37%  ' },
    @{ Cell = 'E20'; Value = '  +0.26%  ' },
    @{ Cell = 'D21'; Value = '72.88' },
    @{ Cell = 'E21'; Value = '  +1.20%  ' },
    @{ Cell = 'D22'; Value = '11.07' },
    @{ Cell = 'E22'; Value = '  -0.40%  ' },
    @{ Cell = 'D23'; Value = '230.55' },
    @{ Cell = 'E23'; Value = '  -0.52%  ' },
    @{ Cell = 'E24'; Value = '  -6.15%  ' },
    @{ Cell = 'E25'; Value = '  +0.02%  ' },
    @{ Cell = 'E26'; Value = '  -3.26%  ' },
    @{ Cell = 'E27'; Value = '  -0.13%  ' },
    @{ Cell = 'D28'; Value = '2.27' },
    @{ Cell = 'E28'; Value = '  -1.23%  ' },
    @{ Cell = 'D29'; Value = '2.20' },
    @{ Cell = 'E29'; Value = '  -2.71%  ' },
    @{ Cell = 'D30'; Value = '166.61' },
    @{ Cell = 'E30'; Value = '  -0.22%  ' },
    @{ Cell = 'D31'; Value = '20.59' },
    @{ Cell = 'E31'; Value = '  -0.47%  ' },
    @{ Cell = 'D32'; Value = '5.64' },
    @{ Cell = 'E32'; Value = '  -3.74%  ' },
    @{ Cell = 'D33'; Value = '0.0801' },
    @{ Cell = 'E33'; Value = '  -1.90%  ' },
    @{ Cell = 'D34'; Value = '29.98' },
    @{ Cell = 'E34'; Value = '  -2.94%  ' },
    @{ Cell = 'E35'; Value = '  -0.50%  ' },
    @{ Cell = 'D36'; Value = '0.110' },
    @{ Cell = 'E36'; Value = '  -7.33%  ' },
    @{ Cell = 'D37'; Value = '4.32' },
    @{ Cell = 'E37'; Value = '  -6.91%  ' },
    @{ Cell = 'D38'; Value = '0.0304' },
    @{ Cell = 'E38'; Value = '  -3.53%  ' },
    @{ Cell = 'D39'; Value = '13.23' },
    @{ Cell = 'E39'; Value = '  -3.80%  ' },
    @{ Cell = 'D40'; Value = '2.13' },
    @{ Cell = 'E40'; Value = '  -2.45%  ' },
    @{ Cell = 'B41'; Value = 'MultiversX' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' },
    @{ Cell = 'D41'; Value = '65.13' },
    @{ Cell = 'E41'; Value = '  +2.36%  ' },
    @{ Cell = 'B42'; Value = 'THORChain' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ Cell = 'D42'; Value = '5.70' },
    @{ Cell = 'E42'; Value = '  -1.12%  ' },
    @{ Cell = 'E43'; Value = '  -1.34%  ' },
    @{ Cell = 'D44'; Value = '8.71' },
    @{ Cell = 'E44'; Value = '  -1.73%  ' },
    @{ Cell = 'D45'; Value = '104.33' },
    @{ Cell = 'E45'; Value = '  -2.11%  ' },
    @{ Cell = 'E46'; Value = '  -1.94%  ' },
    @{ Cell = 'E47'; Value = '  -2.80%  ' },
    @{ Cell = 'E48'; Value = '  -2.48%  ' },
    @{ Cell = 'E49'; Value = '  -0.98%  ' },
    @{ Cell = 'E50'; Value = '  -1.35%  ' },
    @{ Cell = 'D51'; Value = '2.429.02' },
    @{ Cell = 'E51'; Value = '  -1.12%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $col = $u.Cell -replace '[0-9]+$', ''

    if ($col -eq 'D') {
        # Price column: many values look like plain numbers ("0.627",
        # "243.13", ...) but the workbook stores them as text. Force the
        # cell to text format first so Excel doesn't silently reinterpret
        # the assignment as a numeric value, then restore the default
        # "Normal" style so no stray number format sticks around.
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        # Coin name (B), Link (C) and Volume(1h) (E) are never ambiguous
        # with numbers, so a plain assignment is enough.
        $rng.Value = $u.Value
    }
}
